$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "19.903.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -8.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.398.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -8.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3646"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3014"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06431"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9594"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.046"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.29%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.402.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.66%  "
$ws.Range("B16").Value = "Solana"
$ws.Range("C16").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -11.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009914"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05615"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -14.76%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -14.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.482"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.244"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "19.924.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.171"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.560.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "106.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.827"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -20.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.198"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7917"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -15.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07532"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.237"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05647"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.696"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1893"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.925"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.026"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5183"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.457"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.05%  "
$ws.Range("E46").Value = "  -7.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4982"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "108.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.726"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.04%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.023"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.43%  "
